$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..101 down to 9..102
$ws.Rows(8).Insert()

# Populate the new row 8 with the new reading (same market/category block,
# new date, new price tier, new unit, new origin)
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44503
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112031
$ws.Range("G8").Value = "Poroto verde"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 30000
$ws.Range("M8").Value = 30000
$ws.Range("N8").Value = "$/malla 25 kilos"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 1200
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
